$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values C1:S1 (2..18)
$arr1 = New-Object "object[,]" 1,17
$arr1[0,0] = 2
$arr1[0,1] = 3
$arr1[0,2] = 4
$arr1[0,3] = 5
$arr1[0,4] = 6
$arr1[0,5] = 7
$arr1[0,6] = 8
$arr1[0,7] = 9
$arr1[0,8] = 10
$arr1[0,9] = 11
$arr1[0,10] = 12
$arr1[0,11] = 13
$arr1[0,12] = 14
$arr1[0,13] = 15
$arr1[0,14] = 16
$arr1[0,15] = 17
$arr1[0,16] = 18
$ws.Range("C1:S1").Value = $arr1

# Row 2: B2:S2
$arr2 = New-Object "object[,]" 1,18
$arr2[0,0] = -105.887540230357
$arr2[0,1] = -423.7703112008439
$arr2[0,2] = -2257914.998695466
$arr2[0,3] = -796948511.7800217
$arr2[0,4] = -241923537953.5019
$arr2[0,5] = -69556378333898.71
$arr2[0,6] = [double]"-1.929968358403642e+16"
$arr2[0,7] = [double]"-5.214099365283171e+18"
$arr2[0,8] = [double]"-1.379339118921484e+21"
$arr2[0,9] = [double]"-3.570643344586957e+23"
$arr2[0,10] = [double]"-9.055649038096875e+25"
$arr2[0,11] = [double]"-2.261637941677171e+28"
$arr2[0,12] = [double]"-5.602547861823997e+30"
$arr2[0,13] = [double]"-1.381114534473173e+33"
$arr2[0,14] = [double]"-3.391527605264924e+35"
$arr2[0,15] = [double]"-8.310062931270042e+37"
$arr2[0,16] = [double]"-2.035025336630176e+40"
$arr2[0,17] = [double]"-4.984945818600842e+42"
$ws.Range("B2:S2").Value = $arr2

# Row 3: B3:S3
$arr3 = New-Object "object[,]" 1,18
$arr3[0,0] = -91.50156156679026
$arr3[0,1] = -20398.84125187307
$arr3[0,2] = -3729683.643309769
$arr3[0,3] = -708926658.521165
$arr3[0,4] = -137036728453.9774
$arr3[0,5] = -26565758862733.38
$arr3[0,6] = -4933769354479878
$arr3[0,7] = [double]"9.093414429489686e+17"
$arr3[0,8] = [double]"1.648590427360468e+20"
$arr3[0,9] = [double]"2.905324393148314e+22"
$arr3[0,10] = [double]"4.883973608763257e+24"
$arr3[0,11] = [double]"7.573473759333211e+26"
$arr3[0,12] = [double]"1.016384166799409e+29"
$arr3[0,13] = [double]"1.02072447455898e+31"
$arr3[0,14] = [double]"3.515974928385321e+32"
$arr3[0,15] = [double]"-1.298066142634713e+35"
$arr3[0,16] = [double]"-3.606351361404414e+37"
$arr3[0,17] = [double]"-6.20465347621206e+39"
$ws.Range("B3:S3").Value = $arr3

# Row 4: B4:S4
$arr4 = New-Object "object[,]" 1,18
$arr4[0,0] = -37.3856736201392
$arr4[0,1] = 3605.83386231337
$arr4[0,2] = 597608.2784633371
$arr4[0,3] = 76500185.88382883
$arr4[0,4] = 5593945912.290071
$arr4[0,5] = -670998559267.8333
$arr4[0,6] = -180487872261025.8
$arr4[0,7] = [double]"-3.387293773469449e+16"
$arr4[0,8] = [double]"-4.400891470100404e+18"
$arr4[0,9] = [double]"-1.464594869207166e+20"
$arr4[0,10] = [double]"1.251146652340841e+23"
$arr4[0,11] = [double]"4.53638403926917e+25"
$arr4[0,12] = [double]"9.202352703387804e+27"
$arr4[0,13] = [double]"1.297758968214318e+30"
$arr4[0,14] = [double]"2.074088193637579e+32"
$arr4[0,15] = [double]"3.8341445785591e+34"
$arr4[0,16] = [double]"7.449962612914107e+36"
$arr4[0,17] = [double]"1.485561826352686e+39"
$ws.Range("B4:S4").Value = $arr4

# Row 5: B5:S5
$arr5 = New-Object "object[,]" 1,18
$arr5[0,0] = 62.70481496860483
$arr5[0,1] = 6962.007188279983
$arr5[0,2] = 1272930.107868484
$arr5[0,3] = 274696336.7557957
$arr5[0,4] = 57443388700.50078
$arr5[0,5] = 11573909426070.4
$arr5[0,6] = 2140667666749535
$arr5[0,7] = [double]"3.721417228238032e+17"
$arr5[0,8] = [double]"5.88110991290336e+19"
$arr5[0,9] = [double]"-8.001425468955014e+21"
$arr5[0,10] = [double]"-8.724432767196372e+23"
$arr5[0,11] = [double]"-7.460880880823585e+25"
$arr5[0,12] = [double]"-7.930609095198493e+27"
$arr5[0,13] = [double]"-1.378064453508017e+30"
$arr5[0,14] = [double]"-2.407076393761308e+32"
$arr5[0,15] = [double]"-4.292727335298168e+34"
$arr5[0,16] = [double]"-7.863526935818906e+36"
$arr5[0,17] = [double]"-1.460791566723002e+39"
$ws.Range("B5:S5").Value = $arr5

# Row 6: B6:S6
$arr6 = New-Object "object[,]" 1,18
$arr6[0,0] = 46.26634212775927
$arr6[0,1] = 1111.919801866168
$arr6[0,2] = -328307.2507687235
$arr6[0,3] = -38229508.62856968
$arr6[0,4] = -3434500800.793842
$arr6[0,5] = -190710517305.7799
$arr6[0,6] = 52457107918889.59
$arr6[0,7] = [double]"1.385455856756534e+16"
$arr6[0,8] = [double]"1.824401630156151e+18"
$arr6[0,9] = [double]"2.161137216471457e+20"
$arr6[0,10] = [double]"2.640298731810143e+22"
$arr6[0,11] = [double]"3.062068118104962e+24"
$arr6[0,12] = [double]"3.774418677753765e+26"
$arr6[0,13] = [double]"4.779714970973e+28"
$arr6[0,14] = [double]"5.948224143948945e+30"
$arr6[0,15] = [double]"7.961400868661381e+32"
$arr6[0,16] = [double]"1.231759650004992e+35"
$arr6[0,17] = [double]"2.088378438284379e+37"
$ws.Range("B6:S6").Value = $arr6

# Row 7: B7:S7
$arr7 = New-Object "object[,]" 1,18
$arr7[0,0] = -4.280939842960842
$arr7[0,1] = -5253.121963825627
$arr7[0,2] = 651009.3540266006
$arr7[0,3] = 81063511.76994532
$arr7[0,4] = 9158026691.554754
$arr7[0,5] = 823369565815.3275
$arr7[0,6] = -26370054848166.36
$arr7[0,7] = 5527602228218800
$arr7[0,8] = [double]"8.974825376268806e+17"
$arr7[0,9] = [double]"1.143326611989307e+20"
$arr7[0,10] = [double]"1.405368061884534e+22"
$arr7[0,11] = [double]"1.905551828097895e+24"
$arr7[0,12] = [double]"3.297420638626605e+26"
$arr7[0,13] = [double]"6.338896669922746e+28"
$arr7[0,14] = [double]"1.164484694771914e+31"
$arr7[0,15] = [double]"2.02612310914652e+33"
$arr7[0,16] = [double]"3.360077956779846e+35"
$arr7[0,17] = [double]"5.440948609863128e+37"
$ws.Range("B7:S7").Value = $arr7

# Apply header style (bold, border, center) from B1 to the new header cells C1:S1
$ws.Range("B1").Copy()
$ws.Range("C1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A1").Select()
